# Update LR-pairs data for Cort-Sstr2 with recomputed TPM-based values.
# Each cell below is set to the exact numeric value produced by the
# updated scripts (new TPM input), matching the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.651114
$ws.Range("I2").Value = 0.8572432933444277
$ws.Range("J2").Value = 0.8572432933444277
$ws.Range("M2").Value = 0.097952
$ws.Range("N2").Value = 0.293856
$ws.Range("O2").Value = 0.3056806443660103
$ws.Range("P2").Value = 0.3056806443660104
$ws.Range("Q2").Value = 0.021259306176
$ws.Range("R2").Value = 0.191333755584
$ws.Range("S2").Value = 0.2620426822879655
$ws.Range("T2").Value = 0.2620426822879655
$ws.Range("H3").Value = 0.651114
$ws.Range("I3").Value = 0.8572432933444277
$ws.Range("J3").Value = 0.8572432933444277
$ws.Range("O3").Value = 0.3650866467564809
$ws.Range("P3").Value = 0.3650866467564809
$ws.Range("S3").Value = 0.3129680794215994
$ws.Range("T3").Value = 0.3129680794215995
$ws.Range("H4").Value = 0.651114
$ws.Range("I4").Value = 0.8572432933444277
$ws.Range("J4").Value = 0.8572432933444277
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.105499
$ws.Range("N4").Value = 0.316497
$ws.Range("O4").Value = 0.3292327088775087
$ws.Range("P4").Value = 0.3292327088775087
$ws.Range("Q4").Value = 0.022897291962
$ws.Range("R4").Value = 0.206075627658
$ws.Range("S4").Value = 0.2822325316348627
$ws.Range("T4").Value = 0.2822325316348628
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03614333333333333
$ws.Range("H5").Value = 0.10843
$ws.Range("I5").Value = 0.1427567066555723
$ws.Range("J5").Value = 0.1427567066555723
$ws.Range("M5").Value = 0.097952
$ws.Range("N5").Value = 0.293856
$ws.Range("O5").Value = 0.3056806443660103
$ws.Range("P5").Value = 0.3056806443660104
$ws.Range("Q5").Value = 0.003540311786666667
$ws.Range("R5").Value = 0.03186280608
$ws.Range("S5").Value = 0.04363796207804486
$ws.Range("T5").Value = 0.04363796207804486
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03614333333333333
$ws.Range("H6").Value = 0.10843
$ws.Range("I6").Value = 0.1427567066555723
$ws.Range("J6").Value = 0.1427567066555723
$ws.Range("O6").Value = 0.3650866467564809
$ws.Range("P6").Value = 0.3650866467564809
$ws.Range("Q6").Value = 0.00422833628
$ws.Range("R6").Value = 0.03805502652
$ws.Range("S6").Value = 0.05211856733488149
$ws.Range("T6").Value = 0.0521185673348815
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03614333333333333
$ws.Range("H7").Value = 0.10843
$ws.Range("I7").Value = 0.1427567066555723
$ws.Range("J7").Value = 0.1427567066555723
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.105499
$ws.Range("N7").Value = 0.316497
$ws.Range("O7").Value = 0.3292327088775087
$ws.Range("P7").Value = 0.3292327088775087
$ws.Range("Q7").Value = 0.003813085523333334
$ws.Range("R7").Value = 0.03431776971
$ws.Range("S7").Value = 0.04700017724264594
$ws.Range("T7").Value = 0.04700017724264595
